{"js": "// The \"Date Range\" column of the StreamInfo table lists spans such as\n// \"2019-2023\", \"2017-2023\", etc. The survey window's end-year moved\n// from 2023 to 2022 (e.g. \"2019-2023\" -> \"2019-2022\"). Date ranges that\n// already ended in a different year (e.g. \"2018-2022\", \"2017-2021\")\n// must stay untouched, so only search for/replace the specific\n// \"####-2023\" strings that actually appear in the table.\nconst body = context.document.body;\n\nconst dateRanges = [\n  \"2019-2023\",\n  \"2017-2023\",\n  \"2015-2023\",\n  \"2020-2023\",\n  \"2016-2023\"\n];\n\nfor (const oldRange of dateRanges) {\n  const newRange = oldRange.slice(0, 5) + \"2022\";\n  const results = body.search(oldRange, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newRange, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The \"Table 1: StreamInfo\" table lists a \"Date Range\" column such as\n# \"2019-2023\", \"2017-2023\", etc. The survey window end-year was updated\n# from 2023 to 2022 (e.g. \"2019-2023\" -> \"2019-2022\"), while ranges that\n# already ended in a different year (e.g. \"2018-2022\", \"2017-2021\") are\n# left untouched. Do a targeted Find & Replace for every \"####-2023\"\n# date-range string, replacing just the trailing year.\n$d = $word.ActiveDocument\n\n$ranges = @(\n    \"2019-2023\",\n    \"2017-2023\",\n    \"2015-2023\",\n    \"2020-2023\",\n    \"2016-2023\"\n)\n\nforeach ($old in $ranges) {\n    $new = $old.Substring(0, 5) + \"2022\"\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($null, $false, $false, $false, $null, $null, $true, $null, $null, $null, 2)\n}\n"}
